# schedule.xlsx update: mark row 15 ("Android多媒体..." task) as completed on
# 2018-05-22, and append a new schedule row (row 16) for the next task
# ("Android应用资源学习(res目录下)") assigned 2018-05-27.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 15: fill in the completion date (C15) -----------------------------
$ws.Range("C15").Value = 43242
$ws.Range("C15").NumberFormat = "m/d/yy"

# --- Row 16: new task entry -------------------------------------------------
# Copy the formatting (border / wrap-text / font) from the row above so the
# new row matches the rest of the table, then overwrite with the new values.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A16").Value = "1.Android应用资源学习(res目录下)`n"

$ws.Range("B15").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("B16").Value = 43247

$ws.Range("D15").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D16").Value = "疯狂Android讲义第6章(疯狂Android比较杂，以后可以当成`"字典`"来查)，系统学一下res目录下的文件编写及存放目录"

$ws.Rows.Item(16).RowHeight = 135

# --- Selection follows the newly-added row ---------------------------------
[void]$ws.Range("E16").Select()
